$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 15 simulation numbers to the long data lists ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "4 2 11 5 4 33 4 5 4 4 2 13 3 2 0 5 9 1 1 4 6 9 5 4 -3 6 -3 2 7 5 1 0 2 4 4 2 25 1 8 1 -3 -1 5 2 12 2 7 4 2 2 1 8 8 1 11 8 3 3 30 5 6 1 4 2 4 2 -2 -1 5 4 3 5 2 12 -1 6 8 6 16 1 6 4 3 -2 4 14 0 1 5 3 13 1 4 4 0 7 4 2 3 4 4 3 6 -4 5 4 2 2 6 2 1 5 12 7 7 6 -1 5 11 21 5 7 7 5 0 2 0 3 1 14 3 26 5 10 -1 -1 6 1 -4 3 6 0 0 1 8 12 3 -3 1 1 5 12 2 14 3 6 1 4 5 0 1 -1 5 4 -1 -1 12 9 2 5 0 0 0 2 7 5 -5 4 3 -5 2 2 3 7 6 4 2 -2 -2 7 2 7 10 6 0 5 4 5 1 4 7 1 12 5 2 4 18 4 5 2 3 11 4 0 4 9 2 0 1 2 -4 7 15 2 1 -1 2 7 5 0 0 1 3 0 1 12 4 -5 0 11 12 0 7 4 2 1 4 2 3 3 3 3 2 0 4 1 10 4 2 0 -1 9 13 -1 3 -3 5 3 32 7 -3 8 5 9 3 13 -1 3 13 3 0 3 6 3 47 4 3 0 4 2"
$ws.Range("C2").Value = "1 6 7 4 8 9 2 -3 3 2 3 1 1 2 4 3 2 -2 4 9 4 8 3 12 1 1 7 3 8 1 -2 2 7 0 1 4 7 3 2 4 8 4 1 11 3 4 -2 9 2 6 7 0 60 3 6 2 19 0 3 8 4 4 1 0 12 4 1 4 7 1 3 1 4 8 2 2 12 2 2 9 3 4 3 6 3 -1 5 0 2 24 9 2 3 4 13 2 0 1 1 -2 2 3 5 6 9 2 7 4 2 8 4 5 5 -1 5 15 7 -1 1 6 7 2 4 3 3 6 1 7 7 1 4 3 -3 10 2 6 15 3 6 0 1 3 -1 7 5 29 8 -1 0 2 12 0 2 -1 0 4 1 5 0 -2 1 6 1 25 2 4 1 4 4 5 0 1 6 -1 10 11 2 1 4 9 3 10 5 4 -3 5 4 2 3 1 3 4 3 2 0 0 2 18 3 -1 6 2 5 -3 3 9 1 7 1 12 4 -2 4 14 4 1 3 2 4 4 3 4 8 8 1 4 3 7 4 5 3 11 5 -1 3 7 5 8 1 2 0 2 0 5 4 7 2 1 6 2 2 5 4 3 4 2 2 5 2 0 2 2 0 7 0 8 4 2 1 2 2 5 2 0 7 -1 3 0 0 6 1 2 5 4 2 7 6 3 11 4 1 2 6 7 3 8 4 5 1 0 4 3 4 4 3 1 1 0 1 5 7 3 4 4 3 10 10 2 1 6 3 0 0 7 4 1 7 5 17 1 2 0 2 2 1 -1 4 5 0 6 2 7 5 2 9 2 4 0 2 2 5 4 0 2 6 3 1 -1 4 0 1 6 5 -6 2 10 0 2 2 -1 5 2 9 0 -3 0 8 3 5 -4 7 4 3"
$ws.Range("B3").Value = "23 11 6 22 9 16 4 6 69 4 12 5 3 9 30 6 4 15 51 4 9 4 5 16 63 11 8 12 8 1 14 10 11 1 7 68 3 7 3 27 17 5 28 10 10 12 9 5 21 7 39 7 12 11 31 5 14 6 -3 19 7 17 15 -2 5 3 2 28 28 12 1 19 9 2 13 16 6 2 5 25 4 25 14 19 29 9 8 13 5 1 17 8 3 21 23 7 5 14 20 13 -1 2 27 12 -3 12 14 0 9 14 41 1 3 -4 7 7 7 2 11 7 12 21 10 84 5 -4 5 10 10 1 7 12 9 28 0 9 27 1 5 8 5 9 5 8 16 7 28 5 4 8 5 28 5 9 3 12 -6 6 28 1 9 9 7 5 7 18 8 6 -4 4 8 10 8 9 6 9 25 14 36 4 12 16 6 16 48 2 10 3 2 55 6 13 9 39 13 8 2 7 6 1 12 1 4 13 13 4 7 32 -1 4 1 0 8 10 7 12 -1 33 5 9 2 4 11 7 6 7 16 17 5 7 2 12 11 5 2 27 3 0 1 29 5 8 24 25 16 55 7 13 22 14 1 29 6 1 4"
$ws.Range("C3").Value = "10 6 10 7 24 6 10 14 10 8 4 5 7 15 16 8 9 7 8 13 7 1 16 11 19 10 6 14 13 5 17 51 15 15 19 6 15 15 22 12 12 5 15 16 10 14 4 9 12 4 5 6 20 17 8 10 7 6 23 26 7 15 7 4 7 9 16 3 -4 15 5 8 9 -1 21 11 28 4 14 15 2 11 9 8 21 9 18 21 7 15 14 26 7 10 8 8 14 10 12 12 7 8 -1 16 11 9 4 17 17 25 5 9 5 8 6 7 9 68 20 20 13 3 23 1 24 33 13 7 9 2 7 9 12 6 4 12 5 6 5 12 9 6 5 2 1 15 14 23 6 12 4 12 7 12 2 3 13 8 6 23 7 14 5 9 14 6 31 15 -5 6 7 13 8 11 13 6 1 7 5 17 9 10 13 7 7 10 3 9 7 0 7 6 3 6 11 9 8 11 7 4 5 9 10 9 5 8 3 41 9 4 6 12 6 34 3 4 14 8 10 6 4 23 24 4 7 8 5 50 7 4 8 17 5 7 11 1 -5 31 7 1 12 25 2 10 20 6 11 4 21 6 -2 14 7 5 16 5 3 10 9 15 7 3 12 20 5 4 15 2 8 6 0 1 1 12 9 8 14 9 16 10 20 2 9 11 11 6 3 12 6 8 24 8 28 12 4 7 9 8 9 48 3 21 13 29 11 18 14 10 7 16 7 17 13 30 5 10 7 0 15 8 11 15 15 12 6 9 18 7 11 -2 19 10 7 11 9 11 12 4 17 6 18 7 2 6 8"

# --- OFF sheet: updated aggregate totals ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 147
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 36
$ws.Range("G2").Value = 36
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 22
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 116
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = 72
$ws.Range("G3").Value = 30
$ws.Range("I3").Value = 40
$ws.Range("J3").Value = 38
$ws.Range("L3").Value = 208
$ws.Range("M3").Value = 137
$ws.Range("Q3").Value = 384

# --- DEF sheet: updated aggregate totals ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 182
$ws.Range("D2").Value = 9
$ws.Range("F2").Value = 53
$ws.Range("G2").Value = 65
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 28
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 23
$ws.Range("C3").Value = 180
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 104
$ws.Range("G3").Value = 38
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 44
$ws.Range("J3").Value = 47
$ws.Range("L3").Value = 278
$ws.Range("M3").Value = 191
$ws.Range("Q3").Value = 515

# --- ST sheet: updated aggregate totals + appended simulation lists ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 55
$ws.Range("D2").Value = 69
$ws.Range("F2").Value = 250
$ws.Range("G2").Value = 224
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 87
$ws.Range("K2").Value = 82
$ws.Range("B3").Value = 17
$ws.Range("B4").Value = "70 65 69 63 68 51 64 63 65 66 66 61 55 65 66 58 55 65 66 65 66 69 64 39 68 62 61 68 60 59 55 65 62 55 58 66 61 50"
$ws.Range("B5").Value = "22 19 19 14 19 13 19 16 24 16 26 7 21 16 25 20 12 37 19 17 33 20 15 13 25 22 32 26 0 27 12 31 12 12 23 13 24 11"
$ws.Range("B6").Value = "22 7 18 21 33 18 24 22 0 0 34 38 25 17 23 13 21 24 32 23 20 27 27 27 22 23 25 33 26 19 15 18 25"
$ws.Range("D3").Value = "25 49 40 36 35 66 52 40 59 46 47 43 45 50 60 34 55 62 51 50 43 31 51 52 68 57 43 24 47 42 39 65 47 65 51 32 44 50 48 44 34 41 59 46 45 36 35 51 61 53 42 60 50 52 48 55 54 38 54 50 48 50 65 36 54 41 48 48 32"
$ws.Range("D4").Value = "0 7 0 0 0 0 0 0 10 13 9 13 -1 10 10 0 10 12 0 0 0 0 0 11 0 6 0 0 0 0 0 14 9 0 2 0 0 0 0 13 0 0 4 9 0 0 4 0 16 8 0 0 12 0 0 0 0 0 0 8 0 0 22 0 -4 0 2 0 0"
$ws.Range("D5").Value = "0 0 0 0 0 9 22"

# --- TURNS sheet ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("D3").Value = 7

# --- PEN sheet ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 9
$ws.Range("B3").Value = 9
$ws.Range("D3").Value = 7
$ws.Range("D4").Value = 6
